$d = $word.ActiveDocument
$bull = [char]0x2022

# ---------------------------------------------------------------------
# 1. Contact line: merge the spell-checked / split runs into one run and
#    append the extra "github.com/danielmartincraig" + linkedin text.
# ---------------------------------------------------------------------
$contactNew = "(803)389-6750 $bull danielmartincraig@gmail.com $bull github.com/danielmartincraig $bull linkedin.com/danielcraig23"
$d.Content.Find.Execute("github.com/danielmartincraig $bull linkedin.com/danielcraig23", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$okContact = $d.Content.Find.Execute($contactNew, $true, $false, $false, $false, $false, $true, 1, $false, $contactNew, 2)

# ---------------------------------------------------------------------
# 2. "Web Engineering I and II" bullet: collapse the gramErr-wrapped runs
#    into a single plain run.
# ---------------------------------------------------------------------
$webText = "$bull    Web Engineering I and II"
$okWeb = $d.Content.Find.Execute($webText, $true, $false, $false, $false, $false, $true, 1, $false, $webText, 2)

# ---------------------------------------------------------------------
# 3. Remove the old "_GoBack" bookmark (was at the very end, after
#    "Fluent in Spanish").
# ---------------------------------------------------------------------
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

# ---------------------------------------------------------------------
# 4. Insert a new "OBJECTIVE:" paragraph right after the contact-info
#    paragraph (i.e. right before "EDUCATION:").
# ---------------------------------------------------------------------
$eduPara = $d.Paragraphs(4)

function Insert-TextAtEnd($doc, $para, [string]$text) {
    $pos = $para.Range.End - 1
    $rg = $doc.Range($pos, $pos)
    $rg.InsertAfter($text)
}

# 4a. Create an empty paragraph before EDUCATION, styled Heading1, and
#     type the (unformatted / style-default-sized) "OBJECTIVE: " label.
$r = $eduPara.Range
$r.Collapse(1)
$r.InsertParagraphBefore()
$objPara = $d.Paragraphs(4)
$objPara.Style = "Heading1"
Insert-TextAtEnd $d $objPara "OBJECTIVE: "

# 4b. Create a second, temporary paragraph (also Heading1) that holds
#     "Eager to drive back-end solutions at ", then bump its font size to
#     12pt -- this also stamps the paragraph-mark's own rPr with sz=24,
#     which is what we want once the two paragraphs are merged together.
$objPara.Range.InsertParagraphAfter()
$restPara = $d.Paragraphs(5)
$restPara.Style = "Heading1"
Insert-TextAtEnd $d $restPara "Eager to drive back-end solutions at "
$restPara.Range.Font.Size = 12

# 4c. Merge the two paragraphs by deleting the paragraph mark between
#     them -- "OBJECTIVE: " keeps its bare (no rPr) formatting while the
#     rest of the sentence keeps its sz=24 formatting, and the surviving
#     paragraph mark is the sz=24 one from the second paragraph.
$markStart = $objPara.Range.End - 1
$markEnd = $objPara.Range.End
$markRange = $d.Range($markStart, $markEnd)
$markRange.Delete()

# 4d. Append "Instructure" as its own run (12pt), then stamp the new
#     "_GoBack" bookmark right after it, then append the closing phrase
#     as a final run (12pt).
Insert-TextAtEnd $d $objPara "Instructure"
$prefixLen = ("OBJECTIVE: Eager to drive back-end solutions at ").Length
$instrStart = $objPara.Range.Start + $prefixLen
$instrEnd = $instrStart + ("Instructure").Length
$d.Range($instrStart, $instrEnd).Font.Size = 12

$bmRange = $d.Range($instrEnd, $instrEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Insert-TextAtEnd $d $objPara " on a full-time basis"
$tailStart = $instrEnd
$tailEnd = $objPara.Range.End - 1
$d.Range($tailStart, $tailEnd).Font.Size = 12

Write-Output "done"
